# Apply updated cryptocurrency price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.400.67'
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '1.671.34'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  +0.64%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.99'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5333'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.010'
$ws.Range('E7').Value = '  +0.65%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2662'
$ws.Range('E8').Value = '  +1.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06393'
$ws.Range('E9').Value = '  +1.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.91'
$ws.Range('E10').Value = '  +2.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07850'
$ws.Range('E11').Value = '  +0.60%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.535'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').Value = '1.638.15'
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('D14').Value = '1.900.83'
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5616'
$ws.Range('E15').Value = '  +2.19%  '
$ws.Range('D16').Value = '0.0₅8192'
$ws.Range('E16').Value = '  -0.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.23'
$ws.Range('E17').Value = '  +1.24%  '
$ws.Range('D18').Value = '26.437.26'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.010'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.710'
$ws.Range('E20').Value = '  +2.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '199.54'
$ws.Range('E21').Value = '  +4.49%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.32'
$ws.Range('E22').Value = '  +2.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.066'
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.011'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.62'
$ws.Range('E25').Value = '  +1.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1231'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.247'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.25'
$ws.Range('E28').Value = '  +1.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.498'
$ws.Range('E29').Value = '  +1.92%  '
$ws.Range('E30').Value = '  +3.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.285'
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.557'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.307'
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.614'
$ws.Range('E34').Value = '  +0.99%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9684'
$ws.Range('E35').Value = '  +1.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.844'
$ws.Range('E36').Value = '  +1.56%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.438'
$ws.Range('E37').Value = '  +0.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5823'
$ws.Range('E38').Value = '  +1.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01613'
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('D40').Value = '1.078.70'
$ws.Range('E40').Value = '  +3.94%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8666'
$ws.Range('E41').Value = '  +1.69%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.903'
$ws.Range('E42').Value = '  +1.66%  '
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.84'
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('D45').Value = '1.811.16'
$ws.Range('E45').Value = '  +0.83%  '
$ws.Range('E46').Value = '  +3.07%  '
$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.011'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').Value = '0.0₈105'
$ws.Range('E48').Value = '  -5.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4412'
$ws.Range('E49').Value = '  +1.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.007'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05167'
$ws.Range('E51').Value = '  +0.23%  '
